$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new hospital/clinic records collected in the 05-08-2020 update
# (rows 97-107). Cells are written in natural reading order (row by row,
# left to right) so that any newly-created shared strings land at the
# same indices Excel itself would assign.

# Row 97 - Gracias a Dios / Puerto Lempira
$ws.Range("G97").Value = 'Gracias a Dios'
$ws.Range("K97").Value = 'Puerto Lempira '
$ws.Range("S97").Value = 'hospital '
$ws.Range("U97").Value = 'Hospital'
$ws.Range("V97").Value = 15.260996
$ws.Range("W97").Value = -83.778705

# Row 98 - Olancho / Juticalpa
$ws.Range("G98").Value = 'Olancho'
$ws.Range("K98").Value = 'Juticalpa'
$ws.Range("S98").Value = 'hospital'
$ws.Range("U98").Value = 'Hospital San Francisco'
$ws.Range("V98").Value = 14.665319
$ws.Range("W98").Value = -86.2259

# Row 99 - Olancho / Juticalpa
$ws.Range("G99").Value = 'Olancho'
$ws.Range("K99").Value = 'Juticalpa'
$ws.Range("S99").Value = 'hospital'
$ws.Range("U99").Value = 'Hospital Trochez Montalvan'
$ws.Range("V99").Value = 14.671486
$ws.Range("W99").Value = -86.221148

# Row 100 - Olancho / Catacamas
$ws.Range("G100").Value = 'Olancho'
$ws.Range("K100").Value = 'Catacamas'
$ws.Range("S100").Value = 'hospital'
$ws.Range("U100").Value = 'Hospital y Clínica Campos'
$ws.Range("V100").Value = 14.853908
$ws.Range("W100").Value = -85.894443

# Row 101 - Olancho / Catacamas
$ws.Range("G101").Value = 'Olancho'
$ws.Range("K101").Value = 'Catacamas'
$ws.Range("S101").Value = 'clínica'
$ws.Range("U101").Value = 'Clínica Medicentro'
$ws.Range("V101").Value = 14.851229
$ws.Range("W101").Value = -85.895389

# Row 102 - Olancho / Catacamas
$ws.Range("G102").Value = 'Olancho'
$ws.Range("K102").Value = 'Catacamas'
$ws.Range("S102").Value = 'clínica'
$ws.Range("U102").Value = 'Clínica San Lucas'
$ws.Range("V102").Value = 14.850133
$ws.Range("W102").Value = -85.895661

# Row 103 - La Paz / Guajiquiro
$ws.Range("G103").Value = 'La Paz'
$ws.Range("K103").Value = 'Guajiquiro'
$ws.Range("S103").Value = 'centro de salud '
$ws.Range("U103").Value = 'Centro de Salud El Guajiquiro'
$ws.Range("V103").Value = 14.120423
$ws.Range("W103").Value = -87.829441

# Row 104 - La Paz / Marcala
$ws.Range("G104").Value = 'La Paz'
$ws.Range("K104").Value = 'Marcala'
$ws.Range("S104").Value = 'centro de salud '
$ws.Range("U104").Value = 'Centro de Salud Marcala'
$ws.Range("V104").Value = 14.159074
$ws.Range("W104").Value = -88.036271

# Row 105 - La Paz / San José
$ws.Range("G105").Value = 'La Paz'
$ws.Range("K105").Value = 'San José'
$ws.Range("S105").Value = 'centro de salud '
$ws.Range("U105").Value = 'Centro de Salud'
$ws.Range("V105").Value = 14.248234
$ws.Range("W105").Value = -87.959115

# Row 106 - La Paz / La Paz
$ws.Range("G106").Value = 'La Paz'
$ws.Range("K106").Value = 'La Paz'
$ws.Range("S106").Value = 'hospital'
$ws.Range("U106").Value = 'Hospital Montecillos'
$ws.Range("V106").Value = 14.319006
$ws.Range("W106").Value = -87.68094

# Row 107 - Hospital Roberto Suazo Córdova
$ws.Range("S107").Value = 'hospital'
$ws.Range("U107").Value = 'Hospiital Roberto Suazo Córdova'
$ws.Range("V107").Value = 14.32233
$ws.Range("W107").Value = -87.678188

# Expand the HOSPITALES_HN table (and its AutoFilter) to cover the new rows
$lo = $ws.ListObjects.Item("HOSPITALES_HN")
$lo.Resize($ws.Range("A1:W128"))

# Keep the hidden _FilterDatabase name in sync with the new table range
$fdb = $wb.Names.Item("HOSPITALES!_FilterDatabase")
$fdb.RefersTo = '=HOSPITALES!$A$1:$W$128'

# Restore the active selection used while entering the new rows
$ws.Range("L104").Select()
